# datacamp.xlsx: insert a new "Understanding Data Engineering" row above the
# existing row 42 ("Machine Learning with scikit-learn" ...), pushing the
# following rows (old 43-46) down by one (new 44-47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 42 (shifts old rows 42-46 down to 43-47).
$ws.Rows.Item(42).Insert()

# Fill the new row 42 with the course name + its single rating value.
$ws.Range("A42").Value = "Understanding Data Engineering"
$ws.Range("I42").Value = 4

# Match the saved selection/viewport from the edited workbook.
[void]$ws.Range("I43").Select()
